$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the Price (D) cells that will receive updated, numeric-looking text
# values as Text format first, so Excel COM does not silently convert the
# assigned string into a genuine number (these columns are plain text in the
# source data, e.g. "27.434.52" / "321.81" are price strings, not numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# --- Rows where both Price (D) and Volume(1h) (E) change ---
$ws.Range("D2").Value = '27.434.52'
$ws.Range("E2").Value = '  -2.52%  '
$ws.Range("D3").Value = '1.744.76'
$ws.Range("E3").Value = '  -3.07%  '
$ws.Range("D5").Value = '321.81'
$ws.Range("E5").Value = '  -4.15%  '
$ws.Range("D7").Value = '0.4199'
$ws.Range("E7").Value = '  -9.42%  '
$ws.Range("D8").Value = '0.3587'
$ws.Range("E8").Value = '  -3.01%  '
$ws.Range("D9").Value = '45.47'
$ws.Range("E9").Value = '  +0.60%  '
$ws.Range("D10").Value = '0.07418'
$ws.Range("E10").Value = '  -2.43%  '
$ws.Range("D11").Value = '1.114'
$ws.Range("E11").Value = '  -2.85%  '
$ws.Range("D13").Value = '21.48'
$ws.Range("E13").Value = '  -3.81%  '
$ws.Range("D14").Value = '6.101'
$ws.Range("E14").Value = '  -3.62%  '
$ws.Range("D15").Value = '7.194'
$ws.Range("E15").Value = '  -3.43%  '
$ws.Range("D16").Value = '1.736.30'
$ws.Range("E16").Value = '  -3.64%  '
$ws.Range("D17").Value = '0.00001068'
$ws.Range("E17").Value = '  -2.73%  '
$ws.Range("D18").Value = '87.68'
$ws.Range("E18").Value = '  +7.30%  '
$ws.Range("D19").Value = '0.06076'
$ws.Range("E19").Value = '  -9.38%  '
$ws.Range("D21").Value = '16.86'
$ws.Range("E21").Value = '  -3.64%  '
$ws.Range("D22").Value = '6.104'
$ws.Range("E22").Value = '  -4.87%  '
$ws.Range("D23").Value = '0.5225'
$ws.Range("E23").Value = '  -5.55%  '
$ws.Range("D24").Value = '27.459.85'
$ws.Range("E24").Value = '  -2.44%  '
$ws.Range("D25").Value = '11.46'
$ws.Range("E25").Value = '  -3.56%  '
$ws.Range("D26").Value = '2.343'
$ws.Range("E26").Value = '  -2.80%  '
$ws.Range("D27").Value = '20.36'
$ws.Range("E27").Value = '  -1.70%  '
$ws.Range("D30").Value = '1.937.32'
$ws.Range("E30").Value = '  -3.45%  '
$ws.Range("D31").Value = '125.73'
$ws.Range("E31").Value = '  -5.45%  '
$ws.Range("D32").Value = '1.194'
$ws.Range("E32").Value = '  -4.96%  '
$ws.Range("D33").Value = '5.691'
$ws.Range("E33").Value = '  -2.87%  '
$ws.Range("D34").Value = '0.09117'
$ws.Range("E34").Value = '  -4.65%  '
$ws.Range("D35").Value = '3.629'
$ws.Range("E35").Value = '  -9.96%  '
$ws.Range("D36").Value = '12.69'
$ws.Range("E36").Value = '  +5.07%  '
$ws.Range("D37").Value = '0.02294'
$ws.Range("E37").Value = '  -2.42%  '
$ws.Range("D39").Value = '5.078'
$ws.Range("E39").Value = '  -3.25%  '
$ws.Range("D40").Value = '0.06053'
$ws.Range("E40").Value = '  -4.93%  '
$ws.Range("D41").Value = '0.6391'
$ws.Range("E41").Value = '  -3.60%  '
$ws.Range("D43").Value = '1.429'
$ws.Range("E43").Value = '  -5.00%  '
$ws.Range("D46").Value = '13.75'
$ws.Range("E46").Value = '  -3.17%  '
$ws.Range("D47").Value = '3.711'
$ws.Range("E47").Value = '  -2.97%  '
$ws.Range("D48").Value = '0.5850'
$ws.Range("E48").Value = '  -4.21%  '
$ws.Range("D49").Value = '125.54'
$ws.Range("E49").Value = '  -3.44%  '
$ws.Range("D50").Value = '1.949'
$ws.Range("E50").Value = '  -5.04%  '
$ws.Range("D51").Value = '0.06823'
$ws.Range("E51").Value = '  -4.51%  '

# --- Rows where only Volume(1h) (E) changes ---
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("E38").Value = '  -3.90%  '
$ws.Range("E42").Value = '  -3.88%  '

# --- Rows 28/29 and 44/45: Coin name and Link swapped between the pair,
#     with Price/Volume updated to the values for the coin now in that row ---
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '2.385'
$ws.Range("E28").Value = '  +0.19%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = '151.99'
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '7.934'
$ws.Range("E44").Value = '  -1.91%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").Value = '0.9999'
$ws.Range("E45").Value = '  -0.04%  '

# Restore the default (Normal) cell style on the Price cells so only the
# underlying value/type changed -- formatting stays identical to before.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
